$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell A1 was renamed from "judul_proposal" to "judul".
$ws.Range("A1").Value = "judul"

# The active selection moved from D18 to A12.
$ws.Range("A12").Select()
